# Regenerate handoff report: new source file UUID, new xliff hashes, updated timestamps.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

$oldGuid = "63c621f8-4ee2-4fc1-af88-6487fb7cc4ec"
$newGuid = "03e11c4c-a08f-4387-9923-82bd72a53dbf"

$oldHash = "be71c04f89fbb70191f817313170e35adacbfd3e"
$newHash = "dc0d75ec29b062cf252d2f097c0ae1fa7f8445ca"

$newMdName      = "$newGuid.md"
$newMdPath      = "e2e\$newGuid.md"
$newZhCnXlf     = "$newGuid.$newHash.zh-cn.xlf"
$newDeDeXlf     = "$newGuid.$newHash.de-de.xlf"

$newHoDate      = "2016-08-13 07:12:50"
$newZhCnHoDate  = "2016-08-13 07:12:42"

# --- Overview sheet ---
$wsOverview.Range("A2").Value = $newMdName
$wsOverview.Range("B2").Value = $newMdPath
$wsOverview.Range("G2").Value = $newHoDate

# --- zh-cn sheet ---
$wsZhCn.Range("A2").Value = $newMdName
$wsZhCn.Range("G2").Value = $newZhCnXlf
$wsZhCn.Range("H2").Value = $newZhCnHoDate

# --- de-de sheet ---
$wsDeDe.Range("A2").Value = $newMdName
$wsDeDe.Range("G2").Value = $newDeDeXlf
$wsDeDe.Range("H2").Value = $newHoDate

# --- Hyperlink display text updates ---
$wsOverview.Hyperlinks.Item(1).TextToDisplay = $newMdPath
$wsZhCn.Hyperlinks.Item(1).TextToDisplay = $newMdName
$wsDeDe.Hyperlinks.Item(1).TextToDisplay = $newMdName

# --- Column A width normalized to 40 (raw OOXML width) on all three sheets ---
$wsOverview.Columns.Item(1).ColumnWidth = 39.1836734693
$wsZhCn.Columns.Item(1).ColumnWidth = 39.1836734693
$wsDeDe.Columns.Item(1).ColumnWidth = 39.1836734693
